# Updated cryptos list on Tue Oct  8 20:31:26 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto ranking rows, and reflects that ImmutableX / PolygonEcosystemToken
# swapped ranking positions (rows 37 and 38).
#
# D-column prices are digit-grouped / decimal text (e.g. "2.443.37") that
# must stay plain text, not be coerced into a number by Excel's usual
# numeric-literal parsing, so each is written with a leading apostrophe
# (forces text entry, same as typing '62.234.98 into a cell) and the
# resulting quote-prefix style is cleared right back off again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'62.234.98"
$c.Style = "Normal"
$ws.Range("E2").Value = '  -1.53%  '
$c = $ws.Range("D3")
$c.Value = "'2.443.37"
$c.Style = "Normal"
$ws.Range("E3").Value = '  -0.22%  '
$ws.Range("E4").Value = '  -0.08%  '
$c = $ws.Range("D5")
$c.Value = "'583.11"
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.97%  '
$c = $ws.Range("D6")
$c.Value = "'143.90"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -1.22%  '
$c = $ws.Range("D9")
$c.Value = "'2.441.01"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("E10").Value = '  -3.41%  '
$ws.Range("E11").Value = '  +2.53%  '
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("E13").Value = '  -3.18%  '
$ws.Range("E14").Value = '  -2.14%  '
$ws.Range("E15").Value = '  -3.63%  '
$c = $ws.Range("D16")
$c.Value = "'2.864.30"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'62.101.83"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.48%  '
$c = $ws.Range("D18")
$c.Value = "'2.433.10"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -0.20%  '
$c = $ws.Range("D20")
$c.Value = "'7.14"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.91%  '
$c = $ws.Range("D21")
$c.Value = "'330.32"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.57%  '
$c = $ws.Range("D22")
$c.Value = "'4.11"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.54%  '
$ws.Range("E23").Value = '  -4.66%  '
$ws.Range("E24").Value = '  +0.02%  '
$c = $ws.Range("D25")
$c.Value = "'65.77"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.21%  '
$c = $ws.Range("D26")
$c.Value = "'9.39"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +4.64%  '
$c = $ws.Range("D27")
$c.Value = "'620.10"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("E29").Value = '  -6.97%  '
$ws.Range("E30").Value = '  -0.25%  '
$c = $ws.Range("D31")
$c.Value = "'1.44"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -4.44%  '
$c = $ws.Range("D32")
$c.Value = "'8.00"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -2.95%  '
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("E34").Value = '  -0.49%  '
$ws.Range("E35").Value = '  -5.21%  '
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range("D37")
$c.Value = "'0.379"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.25%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D38")
$c.Value = "'1.43"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -6.30%  '
$c = $ws.Range("D39")
$c.Value = "'151.53"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.67%  '
$c = $ws.Range("D40")
$c.Value = "'18.33"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.42%  '
$c = $ws.Range("D41")
$c.Value = "'5.26"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("E42").Value = '  -1.33%  '
$ws.Range("E43").Value = '  +1.29%  '
$ws.Range("E44").Value = '  +0.01%  '
$c = $ws.Range("D45")
$c.Value = "'2.48"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -5.54%  '
$ws.Range("E46").Value = '  -3.65%  '
$c = $ws.Range("D47")
$c.Value = "'3.64"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -3.35%  '
$c = $ws.Range("D48")
$c.Value = "'0.0525"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -1.49%  '
$ws.Range("E49").Value = '  -0.43%  '
$c = $ws.Range("D50")
$c.Value = "'19.51"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -8.04%  '
$ws.Range("E51").Value = '  +6.94%  '